$d = $word.ActiveDocument

# Locate the paragraph that begins the "Chowdhury, A.S., Khaledian, E. ..."
# citation (the one about Capreomycin resistance) and the empty
# ListParagraph immediately following it, then delete both paragraphs in
# their entirety (text + paragraph marks), leaving the preceding
# "Antimicrobial Resistance Prediction..." paragraph and the later empty
# paragraphs untouched.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Khaledian") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target)
    $nextPara = $d.Paragraphs.Item($target + 1)

    $startPos = $startPara.Range.Start
    $endPos = $nextPara.Range.End

    $rng = $d.Range($startPos, $endPos)
    $rng.Delete()
    Write-Output "Deleted paragraphs starting at index $target"
} else {
    Write-Output "Target paragraph not found"
}
